# Regenerate the localization-status report for archive:
#  - flip the "Ready for handoff" status to "In Translation" everywhere it
#    appears (Overview summary columns + each locale sheet's Status column)
#  - the Status column shrinks to fit the new (shorter) text, same as the
#    original report generator would have done via column auto-fit/resize

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-locale status lives in columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Locale sheets: status lives in column C
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Resize the now-narrower Status columns to match the shorter text
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
